# Add three new offense rows (73-75) for case 21TRD09437 / Bunner to Sheet1,
# mirroring rows 70-72 but with a "Guilty" plea instead of "No Contest".
# (Per commit message: related to wiring up a "dismissed" checkbox on the
# NoJail sheet; the data change itself is these appended rows.)

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

function Set-TextCell($row, $col, $text) {
    # Force a cell to hold TEXT even when the value looks numeric
    # (e.g. "4510.11", "0"), matching the source data's string typing,
    # then reset the style so no stray "quote prefix" number format sticks.
    $cell = $ws.Cells.Item($row, $col)
    $cell.Value = "'" + $text
    $cell.Style = "Normal"
}

# Row 73: DUS
$ws.Cells.Item(73, 1).Value = "21TRD09437"
$ws.Cells.Item(73, 2).Value = "Bunner"
$ws.Cells.Item(73, 3).Value = "DUS"
Set-TextCell 73 4 "4510.11"
$ws.Cells.Item(73, 5).Value = "M1"
$ws.Cells.Item(73, 6).Value = "Guilty"
$ws.Cells.Item(73, 7).Value = "Guilty"
$ws.Cells.Item(73, 8).Value = 0
Set-TextCell 73 9 "0"

# Row 74: 1ST SPEED 1 YR SCHOOL >35MPH M4
$ws.Cells.Item(74, 1).Value = "21TRD09437"
$ws.Cells.Item(74, 2).Value = "Bunner"
$ws.Cells.Item(74, 3).Value = "1ST SPEED 1 YR SCHOOL >35MPHM4"
$ws.Cells.Item(74, 4).Value = "4511.21B1A"
$ws.Cells.Item(74, 5).Value = "M4"
$ws.Cells.Item(74, 6).Value = "Guilty"
$ws.Cells.Item(74, 7).Value = "Guilty"
$ws.Cells.Item(74, 8).Value = 0
Set-TextCell 74 9 "0"

# Row 75: RECKLESS OPERATION 1ST IN 1 YR
$ws.Cells.Item(75, 1).Value = "21TRD09437"
$ws.Cells.Item(75, 2).Value = "Bunner"
$ws.Cells.Item(75, 3).Value = "RECKLESS OPERATION 1ST IN 1 YR"
Set-TextCell 75 4 "4511.20"
$ws.Cells.Item(75, 5).Value = "MM"
$ws.Cells.Item(75, 6).Value = "Guilty"
$ws.Cells.Item(75, 7).Value = "Guilty"
$ws.Cells.Item(75, 8).Value = 0
Set-TextCell 75 9 "0"
